$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6551
$ws1.Range("F4").Value = 103
$ws1.Range("F5").Value = 146
$ws1.Range("F7").Value = 78
$ws1.Range("F8").Value = 583

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6552
$ws4.Range("F5").Value = 103
$ws4.Range("F6").Value = 146
$ws4.Range("F9").Value = 78
$ws4.Range("F10").Value = 583
